$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to text format so numeric-looking strings
# (e.g. "1.008", "23.80") are preserved exactly as text rather than
# being auto-converted into numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '29.510.57'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '1.922.42'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = '325.61'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').Value = '0.4840'
$ws.Range('E7').Value = '  +3.34%  '
$ws.Range('D8').Value = '0.4098'
$ws.Range('E8').Value = '  +2.30%  '
$ws.Range('D9').Value = '0.08186'
$ws.Range('E9').Value = '  +2.62%  '
$ws.Range('D10').Value = '1.024'
$ws.Range('E10').Value = '  +3.47%  '
$ws.Range('D11').Value = '23.80'
$ws.Range('E11').Value = '  +6.68%  '
$ws.Range('D12').Value = '1.890.75'
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').Value = '6.053'
$ws.Range('E13').Value = '  +3.69%  '
$ws.Range('D14').Value = '7.239'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').Value = '91.45'
$ws.Range('E15').Value = '  +3.29%  '
$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').Value = '1.008'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.06771'
$ws.Range('E17').Value = '  +2.38%  '
$ws.Range('E18').Value = '  +1.67%  '
$ws.Range('D19').Value = '17.82'
$ws.Range('E19').Value = '  +2.55%  '
$ws.Range('D20').Value = '1.006'
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('D21').Value = '29.543.19'
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('D22').Value = '5.636'
$ws.Range('E22').Value = '  +2.92%  '
$ws.Range('D23').Value = '11.77'
$ws.Range('E23').Value = '  +1.67%  '
$ws.Range('D24').Value = '2.186'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').Value = '2.149.63'
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('D26').Value = '6.733'
$ws.Range('E26').Value = '  +11.58%  '
$ws.Range('D27').Value = '156.98'
$ws.Range('E27').Value = '  +1.22%  '
$ws.Range('D28').Value = '20.12'
$ws.Range('E28').Value = '  +2.95%  '
$ws.Range('D29').Value = '2.127'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('D30').Value = '120.63'
$ws.Range('E30').Value = '  +3.14%  '
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').Value = '0.09589'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').Value = '5.535'
$ws.Range('E33').Value = '  +3.97%  '
$ws.Range('D34').Value = '3.569'
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('D35').Value = '1.393'
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('D36').Value = '0.02295'
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('E37').Value = '  +1.81%  '
$ws.Range('D38').Value = '1.181'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('D39').Value = '0.5987'
$ws.Range('E39').Value = '  +3.44%  '
$ws.Range('D40').Value = '8.055'
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').Value = '10.82'
$ws.Range('E41').Value = '  +8.21%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').Value = '1.006'
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '0.1866'
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '2.441'
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '1.279'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('D46').Value = '12.51'
$ws.Range('E46').Value = '  +3.40%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.07617'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = '0.5599'
$ws.Range('E48').Value = '  +2.70%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '1.965'
$ws.Range('E49').Value = '  +3.83%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '117.17'
$ws.Range('E50').Value = '  +3.44%  '
$ws.Range('B51').Value = 'MXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D51').Value = '2.443'
$ws.Range('E51').Value = '  +4.92%  '

# Restore default (Normal) style so no stray number-format styles remain
$dataRange.Style = "Normal"
